$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Capital_Investment")

$ws.Range("B5").Value = 19239
$ws.Range("C5").Value = "Average monthly costs (rent, service charge, insurance, business rates)"

$ws.Range("B6").Value = 32

$ws.Range("B7").Value = 218498.6

$ws.Range("B8").Value = "December 2025 Total"
